# ----------------------------------------------------------------------------
# Endring av filstruktur. Ny input fra S-omradet.
#
# Adds a new worksheet "faktisk_barnetillegg" as the first (leftmost) tab,
# containing a gr/ar/utb_bt table (A1:C29), and leaves the three existing
# sheets (mottaker, arlig_mottaker, husholdning) in place after it.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# New sheet is inserted before the first existing sheet by default, which is
# exactly where "faktisk_barnetillegg" ends up relative to mottaker /
# arlig_mottaker / husholdning.
$new = $wb.Worksheets.Add()
$new.Name = "faktisk_barnetillegg"

# Header row
$new.Cells.Item(1,1).Value = "gr"
$new.Cells.Item(1,2).Value = "ar"
$new.Cells.Item(1,3).Value = "utb_bt"
$new.Range("A1:C1").Font.Bold = $true
$new.Range("A1:C1").HorizontalAlignment = -4108

$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = 2015
$new.Cells.Item(2,3).Value = 0.59227154983541108
$new.Cells.Item(3,1).Value = 0
$new.Cells.Item(3,2).Value = 2016
$new.Cells.Item(3,3).Value = 0.53231286907565645
$new.Cells.Item(4,1).Value = 0
$new.Cells.Item(4,2).Value = 2017
$new.Cells.Item(4,3).Value = 0.48333099586148459
$new.Cells.Item(5,1).Value = 0
$new.Cells.Item(5,2).Value = 2018
$new.Cells.Item(5,3).Value = 0.453407001181935
$new.Cells.Item(6,1).Value = 0
$new.Cells.Item(6,2).Value = 2019
$new.Cells.Item(6,3).Value = 0.41848403587113631
$new.Cells.Item(7,1).Value = 0
$new.Cells.Item(7,2).Value = 2020
$new.Cells.Item(7,3).Value = 0.39518488139498231
$new.Cells.Item(8,1).Value = 0
$new.Cells.Item(8,2).Value = 2021
$new.Cells.Item(8,3).Value = 0.38429605926228888
$new.Cells.Item(9,1).Value = 1
$new.Cells.Item(9,2).Value = 2015
$new.Cells.Item(9,3).Value = 0.86040506366128466
$new.Cells.Item(10,1).Value = 1
$new.Cells.Item(10,2).Value = 2016
$new.Cells.Item(10,3).Value = 0.74196127731433459
$new.Cells.Item(11,1).Value = 1
$new.Cells.Item(11,2).Value = 2017
$new.Cells.Item(11,3).Value = 0.6420569123369767
$new.Cells.Item(12,1).Value = 1
$new.Cells.Item(12,2).Value = 2018
$new.Cells.Item(12,3).Value = 0.56977996118388696
$new.Cells.Item(13,1).Value = 1
$new.Cells.Item(13,2).Value = 2019
$new.Cells.Item(13,3).Value = 0.52653972498485035
$new.Cells.Item(14,1).Value = 1
$new.Cells.Item(14,2).Value = 2020
$new.Cells.Item(14,3).Value = 0.49490989356469051
$new.Cells.Item(15,1).Value = 1
$new.Cells.Item(15,2).Value = 2021
$new.Cells.Item(15,3).Value = 0.47273156774746228
$new.Cells.Item(16,1).Value = 0
$new.Cells.Item(16,2).Value = 2015
$new.Cells.Item(16,3).Value = 0.64931920528611098
$new.Cells.Item(17,1).Value = 0
$new.Cells.Item(17,2).Value = 2016
$new.Cells.Item(17,3).Value = 0.59468141685458009
$new.Cells.Item(18,1).Value = 0
$new.Cells.Item(18,2).Value = 2017
$new.Cells.Item(18,3).Value = 0.54609167923204638
$new.Cells.Item(19,1).Value = 0
$new.Cells.Item(19,2).Value = 2018
$new.Cells.Item(19,3).Value = 0.51882744978138029
$new.Cells.Item(20,1).Value = 0
$new.Cells.Item(20,2).Value = 2019
$new.Cells.Item(20,3).Value = 0.48158147473268742
$new.Cells.Item(21,1).Value = 0
$new.Cells.Item(21,2).Value = 2020
$new.Cells.Item(21,3).Value = 0.45548972685248229
$new.Cells.Item(22,1).Value = 0
$new.Cells.Item(22,2).Value = 2021
$new.Cells.Item(22,3).Value = 0.44560559738939681
$new.Cells.Item(23,1).Value = 1
$new.Cells.Item(23,2).Value = 2015
$new.Cells.Item(23,3).Value = 0.86040506366128466
$new.Cells.Item(24,1).Value = 1
$new.Cells.Item(24,2).Value = 2016
$new.Cells.Item(24,3).Value = 0.74196127731433459
$new.Cells.Item(25,1).Value = 1
$new.Cells.Item(25,2).Value = 2017
$new.Cells.Item(25,3).Value = 0.6420569123369767
$new.Cells.Item(26,1).Value = 1
$new.Cells.Item(26,2).Value = 2018
$new.Cells.Item(26,3).Value = 0.56977996118388696
$new.Cells.Item(27,1).Value = 1
$new.Cells.Item(27,2).Value = 2019
$new.Cells.Item(27,3).Value = 0.52653972498485035
$new.Cells.Item(28,1).Value = 1
$new.Cells.Item(28,2).Value = 2020
$new.Cells.Item(28,3).Value = 0.49490989356469051
$new.Cells.Item(29,1).Value = 1
$new.Cells.Item(29,2).Value = 2021
$new.Cells.Item(29,3).Value = 0.47273156774746228

# Leave the selection where the author left it before saving.
[void]$new.Range("D23").Select()
